# MOSFERATU JLCPCB BOM update
#
# The BOM lists one CAPACITOR, 1nF line with designators "C21, C34".
# This edit:
#   1. Splits that line: C34 is reassigned away, C3 is added, leaving
#      the 1nF row with designators "C3, C21".
#   2. Adds a new CAPACITOR, 220nF row (designator "C2").
#   3. Adds a new CAPACITOR, 33nF row (designator "C34" - the part
#      removed from the 1nF line above).
#
# Table1 (A1:C30 -> A1:C32) is sorted by the Comment column, so the two
# brand-new rows are inserted directly at their correctly-sorted
# position (immediately after "CAPACITOR, 1nF" and immediately after
# "CAPACITOR, 330nF" respectively) rather than appended + re-sorted -
# that keeps every other, unrelated row exactly where it already was.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-point the existing "CAPACITOR, 1nF" designator list -------------
# Old: "C21, C34"  -> New: "C3, C21"
# (Leading "'" forces text so Excel doesn't try to interpret the value.)
$ws.Range("B8").Value = "'C3, C21"

# --- 2. Insert "CAPACITOR, 220nF" / C2 / 0603 right after row 8 ------------
$ws.Rows("9:9").Insert()
$ws.Range("A9").Value = "'CAPACITOR, 220nF"
$ws.Range("B9").Value = "'C2"
$ws.Range("C9").Value = "'0603"

# --- 3. Insert "CAPACITOR, 33nF" / C34 / 0603 right after "CAPACITOR, 330nF"
#        (that row is now at row 11 after the insert above) ----------------
$ws.Rows("12:12").Insert()
$ws.Range("A12").Value = "'CAPACITOR, 33nF"
$ws.Range("B12").Value = "'C34"
$ws.Range("C12").Value = "'0603"

# --- 4. Grow the table / autofilter to cover the two new rows --------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:C32"))

# --- 5. Match the saved selection (cell B8) --------------------------------
$ws.Range("B8").Select()
